$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 46040
$ws.Range("B10").Value = 48994
$ws.Range("A10:B10").NumberFormat = $ws.Range("A9:B9").NumberFormat
